# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Offense (OFF) sheet - Road ("R") row totals after the Wild Card game
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 424
$wsOff.Range("C3").Value = 308
$wsOff.Range("D3").Value = 129
$wsOff.Range("E3").Value = 56
$wsOff.Range("F3").Value = 8

# Defense (DEF) sheet - Road ("R") row totals after the Wild Card game
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 500
$wsDef.Range("C3").Value = 352
$wsDef.Range("D3").Value = 105
$wsDef.Range("E3").Value = 50
